# Insert a new weekly record row into the dataset.
# The sheet is a flat table (row 1 = headers, rows 2.. = data) ordered so
# that new observations are inserted right before the current row 625,
# pushing all later rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 625..end down by inserting a new row at 625.
$ws.Rows.Item(625).Insert()

# Populate the newly inserted row 625 with the new weekly observation.
$ws.Cells.Item(625, 1).Value = 6
$ws.Cells.Item(625, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(625, 3).Value = "Metropolitana"
$ws.Cells.Item(625, 4).Value = 45265
$ws.Cells.Item(625, 5).Value = 13
$ws.Cells.Item(625, 6).Value = 100112032
$ws.Cells.Item(625, 7).Value = "Zapallo italiano"
$ws.Cells.Item(625, 8).Value = "Sin especificar"
$ws.Cells.Item(625, 9).Value = "Primera"
$ws.Cells.Item(625, 10).Value = 1000
$ws.Cells.Item(625, 11).Value = 9000
$ws.Cells.Item(625, 12).Value = 10000
$ws.Cells.Item(625, 13).Value = 9600
$ws.Cells.Item(625, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(625, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(625, 16).Value = 192
$ws.Cells.Item(625, 17).Value = 50
$ws.Cells.Item(625, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(625, 4).NumberFormat = $ws.Cells.Item(626, 4).NumberFormat
